# Apply data fixes to the "comp_quantity_inst1" sheet and
# switch the active sheet/selection to match the target workbook state.

$wb = $excel.ActiveWorkbook

# --- Update data on "comp_quantity_inst1" ---
$wsComp = $wb.Worksheets.Item("comp_quantity_inst1")

# Row 2: T1 -> T5, quantity 11 -> 21
$wsComp.Range("C2").Value = 21

# Row 3: T2 -> F1 (was T2 -> T1), quantity 22 -> 11
$wsComp.Range("B3").Value = "F1"
$wsComp.Range("C3").Value = 11

# Row 4: F1 -> F2 (was T5 -> T4), quantity 40 -> 26
$wsComp.Range("A4").Value = "F1"
$wsComp.Range("B4").Value = "F2"
$wsComp.Range("C4").Value = 26

# --- Update selections on each sheet ---
$wsParams = $wb.Worksheets.Item("parameters")
$wsParams.Activate()
$wsParams.Range("B4").Select()

# --- Make "comp_quantity_inst1" the active sheet/tab, with its own selection ---
$wsComp.Activate()
$wsComp.Range("C2").Select()
